$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.333228349685669
$ws.Range("B1").Value = 1.419899940490723
$ws.Range("C1").Value = 1.192110419273376
$ws.Range("D1").Value = 2.24606728553772
$ws.Range("E1").Value = 5.249097347259521
